$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (rows unaffected by the row-shift at the bottom) ---

# Row 3 - evento 113
$ws.Range("C3").Value = 4
$ws.Range("E3").Value = 0.02

# Row 4 - evento 115
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

# Row 5 - evento 155
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 0.13

# Row 6 - evento 210
$ws.Range("D6").Value = 92

# Row 7 - evento 215
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0.06

# Row 11 - evento 300
$ws.Range("C11").Value = 48

# Row 12 - evento 330
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 0.37

# Row 14 - evento 342
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 0.05

# Row 15 - evento 346
$ws.Range("C15").Value = 42

# Row 16 - evento 348
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 0.02

# Row 18 - evento 355
$ws.Range("D18").Value = 1

# Row 19 - evento 356
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0.12

# Row 20 - evento 357
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0.37

# Row 21 - evento 365
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0

# Row 24 - evento 455
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0.27

# Row 26 - evento 535
$ws.Range("C26").Value = 0
$ws.Range("E26").Value = 1

# Row 27 - evento 549
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 0.16

# --- Insert a new row before row 36 so a new disease entry can be added,
#     shifting the old row 36 (Zika) down to row 37 ---
$ws.Rows.Item(36).Insert()

# Re-write rows 31-37 with the shifted / updated disease data
# (the "evento" codes in column A are text, like the rest of the sheet,
#  so the column is formatted as Text to keep the numeric-looking codes
#  from being stored as numbers)
$ws.Range("A31:A37").NumberFormat = "@"

# Row 31 - evento 720 (new entry, replacing old 740 data position)
$ws.Range("A31").Value = "720"
$ws.Range("B31").Value = "Sindrome de rubeola congenita"
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 1

# Row 32 - evento 740
$ws.Range("A32").Value = "740"
$ws.Range("B32").Value = "Sifilis congenita"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 0

# Row 33 - evento 750
$ws.Range("A33").Value = "750"
$ws.Range("B33").Value = "Sifilis gestacional"
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 0.27

# Row 34 - evento 813
$ws.Range("A34").Value = "813"
$ws.Range("B34").Value = "Tuberculosis"
$ws.Range("C34").Value = 7
$ws.Range("D34").Value = 4
$ws.Range("E34").Value = 0.09

# Row 35 - evento 831
$ws.Range("A35").Value = "831"
$ws.Range("B35").Value = "Varicela individual"
$ws.Range("C35").Value = 9
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 0

# Row 36 - evento 850 (new row, inserted)
$ws.Range("A36").Value = "850"
$ws.Range("B36").Value = "Vih/sida/mortalidad por sida"
$ws.Range("C36").Value = 7
$ws.Range("D36").Value = 5
$ws.Range("E36").Value = 0.13

# Row 37 - evento 895 (shifted down from row 36, values unchanged)
$ws.Range("A37").Value = "895"
$ws.Range("B37").Value = "Zika"
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 1
